$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "14/04/2022"
$ws.Range("B12").Value = "13:36"
$ws.Range("C12").Value = "13:44"
$ws.Range("D12").Value = 8
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = "EURUSD"
$ws.Range("J12").Value = -0.08999999999999986
$ws.Range("K12").Value = "LOSS"
